$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = ""

$ws.Range("H21").Value = 17
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""

$ws.Range("H23").Value = 17
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = ""
$ws.Range("N76").Value = ""

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = ""
$ws.Range("N79").Value = ""

$ws.Range("H98").Value = 3930.4473
$ws.Range("I98").Value = 3212.6553
$ws.Range("K98").Value = 3212.6553
$ws.Range("M98").Value = -1714.6553

$ws.Range("H103").Value = 1946.5
$ws.Range("I103").Value = 1946.5
$ws.Range("K103").Value = 5839.5
$ws.Range("M103").Value = -5253.5

$ws.Range("H122").Value = 3930.4473
$ws.Range("I122").Value = 3212.6553
$ws.Range("K122").Value = 9637.965899999999
$ws.Range("M122").Value = -7187.965899999999

$ws.Range("H137").Value = 1743.619
$ws.Range("I137").Value = 1780.8
$ws.Range("K137").Value = 5342.4
$ws.Range("M137").Value = -2792.4

$ws.Range("H138").Value = 5763.132
$ws.Range("J138").Value = 10424.407
$ws.Range("L138").Value = 31273.221
$ws.Range("N138").Value = -41553.221

$ws.Range("H140").Value = 140000
$ws.Range("J140").Value = 150000
$ws.Range("L140").Value = 150000
$ws.Range("N140").Value = -160360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 22729906
$ws.Range("I61").Value = 29414420
$ws.Range("K61").Value = 29414420
$ws.Range("M61").Value = -29414208

$ws.Range("H102").Value = 137183.8
$ws.Range("I102").Value = 203225.9
$ws.Range("K102").Value = 203225.9
$ws.Range("M102").Value = -201603.9

$ws.Range("H122").Value = 15153254
$ws.Range("I122").Value = 1711.2222
$ws.Range("K122").Value = 5133.6666
$ws.Range("M122").Value = -2683.6666

$ws.Range("H136").Value = 22729906
$ws.Range("I136").Value = 29414420
$ws.Range("K136").Value = 88243260
$ws.Range("M136").Value = -88240710

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4053.7778
$ws.Range("I20").Value = 3398
$ws.Range("K20").Value = 3398
$ws.Range("M20").Value = -3151

$ws.Range("H82").Value = 17337.1
$ws.Range("I82").Value = 4716.857
$ws.Range("J82").Value = 46784.332
$ws.Range("K82").Value = 4716.857
$ws.Range("L82").Value = 46784.332
$ws.Range("M82").Value = -4333.857
$ws.Range("N82").Value = -47550.332

$ws.Range("H85").Value = 17337.1
$ws.Range("I85").Value = 4716.857
$ws.Range("J85").Value = 46784.332
$ws.Range("K85").Value = 4716.857
$ws.Range("L85").Value = 46784.332
$ws.Range("M85").Value = -3390.857
$ws.Range("N85").Value = -49436.332

$ws.Range("H86").Value = 13071.583
$ws.Range("I86").Value = 23130.6
$ws.Range("J86").Value = 5886.5713
$ws.Range("K86").Value = 23130.6
$ws.Range("L86").Value = 5886.5713
$ws.Range("M86").Value = -22007.6
$ws.Range("N86").Value = -8132.5713

$ws.Range("H89").Value = 13071.583
$ws.Range("I89").Value = 23130.6
$ws.Range("J89").Value = 5886.5713
$ws.Range("K89").Value = 115653
$ws.Range("L89").Value = 29432.8565
$ws.Range("M89").Value = -110037
$ws.Range("N89").Value = -40664.85649999999

$ws.Range("H99").Value = 3467.0667
$ws.Range("I99").Value = 2500.7273
$ws.Range("J99").Value = 6124.5
$ws.Range("K99").Value = 2500.7273
$ws.Range("L99").Value = 6124.5
$ws.Range("M99").Value = -1002.7273
$ws.Range("N99").Value = -9120.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10004747
$ws.Range("I31").Value = 4020.3333
$ws.Range("J31").Value = 15630156
$ws.Range("K31").Value = 4020.3333
$ws.Range("L31").Value = 15630156
$ws.Range("M31").Value = -3725.3333
$ws.Range("N31").Value = -15630746

$ws.Range("H34").Value = 10004747
$ws.Range("I34").Value = 4020.3333
$ws.Range("J34").Value = 15630156
$ws.Range("K34").Value = 4020.3333
$ws.Range("L34").Value = 15630156
$ws.Range("M34").Value = -3818.3333
$ws.Range("N34").Value = -15630560

$ws.Range("H99").Value = 8470.25
$ws.Range("I99").Value = 5349.5
$ws.Range("K99").Value = 5349.5
$ws.Range("M99").Value = -3851.5

$ws.Range("H107").Value = 1128
$ws.Range("I107").Value = 611
$ws.Range("K107").Value = 611
$ws.Range("M107").Value = 1309

$ws.Range("H122").Value = 2772395.8
$ws.Range("I122").Value = 2315.6155
$ws.Range("K122").Value = 6946.8465
$ws.Range("M122").Value = -4496.8465

$ws.Range("H126").Value = 8470.25
$ws.Range("I126").Value = 5349.5
$ws.Range("K126").Value = 16048.5
$ws.Range("M126").Value = -13578.5

$ws.Range("H132").Value = 86045.96000000001
$ws.Range("J132").Value = 3185.3333
$ws.Range("L132").Value = 9555.999899999999
$ws.Range("N132").Value = -14615.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2421.6
$ws.Range("I5").Value = 384
$ws.Range("J5").Value = 3780
$ws.Range("K5").Value = 1152
$ws.Range("L5").Value = 11340
$ws.Range("M5").Value = -1040
$ws.Range("N5").Value = -11564

$ws.Range("H12").Value = 585.6667
$ws.Range("J12").Value = 750.5
$ws.Range("L12").Value = 2251.5
$ws.Range("N12").Value = -2597.5

$ws.Range("H14").Value = 25281.25
$ws.Range("I14").Value = 25281.25
$ws.Range("K14").Value = 75843.75
$ws.Range("M14").Value = -75670.75

$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""

$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""

$ws.Range("H107").Value = 1092.4
$ws.Range("J107").Value = 1636.5555
$ws.Range("L107").Value = 4909.666499999999
$ws.Range("N107").Value = -8749.666499999999

$ws.Range("H135").Value = 2421.6
$ws.Range("I135").Value = 384
$ws.Range("J135").Value = 3780
$ws.Range("K135").Value = 3456
$ws.Range("L135").Value = 34020
$ws.Range("M135").Value = -921
$ws.Range("N135").Value = -39090

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 38466252
$ws.Range("I122").Value = 5034.222
$ws.Range("J122").Value = 125003990
$ws.Range("K122").Value = 15102.666
$ws.Range("L122").Value = 375011970
$ws.Range("M122").Value = -12652.666
$ws.Range("N122").Value = -375016870

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1979.8572
$ws.Range("I22").Value = 1866.2307
$ws.Range("J22").Value = 2164.5
$ws.Range("K22").Value = 1866.2307
$ws.Range("L22").Value = 2164.5
$ws.Range("M22").Value = -1571.2307
$ws.Range("N22").Value = -2754.5

$ws.Range("H27").Value = 1979.8572
$ws.Range("I27").Value = 1866.2307
$ws.Range("J27").Value = 2164.5
$ws.Range("K27").Value = 1866.2307
$ws.Range("L27").Value = 2164.5
$ws.Range("M27").Value = -1759.2307
$ws.Range("N27").Value = -2378.5

$ws.Range("H46").Value = 1806.2963
$ws.Range("J46").Value = 6500
$ws.Range("L46").Value = 6500
$ws.Range("N46").Value = -6876

$ws.Range("H107").Value = 12822
$ws.Range("I107").Value = 12822
$ws.Range("K107").Value = 12822
$ws.Range("M107").Value = -10902

$ws.Range("H122").Value = 3791979.8
$ws.Range("I122").Value = 3937.04
$ws.Range("K122").Value = 11811.12
$ws.Range("M122").Value = -9361.119999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 20000
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 20000
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 20000
$ws.Range("M41").Value = ""
$ws.Range("N41").Value = -20780

$ws.Range("H100").Value = 77694250
$ws.Range("I100").Value = 101002136
$ws.Range("K100").Value = 202004272
$ws.Range("M100").Value = -202003731

$ws.Range("H107").Value = 2917.7273
$ws.Range("J107").Value = 3186.875
$ws.Range("L107").Value = 9560.625
$ws.Range("N107").Value = -13400.625

$ws.Range("H122").Value = 16671032
$ws.Range("I122").Value = 2209.111
$ws.Range("J122").Value = 66677504
$ws.Range("K122").Value = 6627.333
$ws.Range("L122").Value = 200032512
$ws.Range("M122").Value = -4177.333
$ws.Range("N122").Value = -200037412

$ws.Range("H132").Value = 1825.125
$ws.Range("I132").Value = 2057.2856
$ws.Range("K132").Value = 6171.8568
$ws.Range("M132").Value = -3641.8568
